$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 71, shifting the existing rows 71-81 down to 72-82
$ws.Rows("71:71").Insert()

# Populate the newly inserted row 71 with the new weekly data point
$ws.Range("A71").Value = 9
$ws.Range("B71").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C71").Value = "Metropolitana"
$ws.Range("D71").NumberFormat = $ws.Range("D70").NumberFormat
$ws.Range("D71").Value = 45142
$ws.Range("E71").Value = 13
$ws.Range("F71").Value = 100112010
$ws.Range("G71").Value = "Achicoria"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 90
$ws.Range("K71").Value = 7000
$ws.Range("L71").Value = 7000
$ws.Range("M71").Value = 7000
$ws.Range("N71").Value = "`$/caja 16 unidades"
$ws.Range("O71").Value = "Provincia de Quillota"
$ws.Range("P71").Value = 438
$ws.Range("Q71").Value = 16
$ws.Range("R71").Value = "Hortaliza"
